$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column C ("Förändrad") changes from 45184 to 45186 for every data row (2-232).
$ws.Range("C2:C232").Value = 45186

# 2. Rows 2-8 have their HYPERLINK formulas (columns S,T,V,W,X,Y) gain a
#    second "friendly name" argument equal to the case's designation (col A).
$cols  = @("S","T","V","W","X","Y")
$paths = @("artfynd","kartor","klagomål","klagomålsmail","tillsyn","tillsynsmail")
$exts  = @(".xlsx",".png",".docx",".docx",".docx",".docx")

for ($row = 2; $row -le 8; $row++) {
    $name = $ws.Cells.Item($row, 1).Value2
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $url = "https://klasma.github.io/Logging_KUNGSBACKA/" + $paths[$j] + "/" + $name + $exts[$j]
        $formula = '=HYPERLINK("' + $url + '", "' + $name + '")'
        $ws.Range($cols[$j] + $row).Formula = $formula
    }
}

# 3. Row 232 gains an explicit custom row height (15pt).
$ws.Rows.Item(232).RowHeight = 15

# 4. A brand-new row 233 is appended with a new case.
$ws.Range("A233").Value = "A 43526-2023"
$ws.Range("B233").Value = 45184
$ws.Range("B233").NumberFormat = "YYYY-MM-DD"
$ws.Range("C233").Value = 45186
$ws.Range("C233").NumberFormat = "YYYY-MM-DD"
$ws.Range("D233").Value = "HALLANDS LÄN"
$ws.Range("E233").Value = "KUNGSBACKA"
$ws.Range("G233").Value = 1.2
$ws.Range("H233").Value = 0
$ws.Range("I233").Value = 0
$ws.Range("J233").Value = 0
$ws.Range("K233").Value = 0
$ws.Range("L233").Value = 0
$ws.Range("M233").Value = 0
$ws.Range("N233").Value = 0
$ws.Range("O233").Value = 0
$ws.Range("P233").Value = 0
$ws.Range("Q233").Value = 0
$ws.Range("R233").WrapText = $true
